$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction: B27 was a mis-keyed value missing its decimal point ---
$ws.Range("B27").Value = 40.44194444

# Row 27's height reverts to the sheet default (12.8) instead of the old 13.4 override
$ws.Rows.Item(27).RowHeight = 12.8

# --- Column widths: give A/B/C their own (now-differentiated) widths ---
# (ColumnWidth is expressed in character units and gets quantized to the
# nearest pixel by the engine, so these are chosen to land as close as
# possible to the authored widths of 11.34 / 15.64 / 17.40.)
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334
$ws.Columns.Item(3).ColumnWidth = 16.5

# --- View: scroll / reselect to reflect where the editor ended up working ---
$win = $excel.ActiveWindow
$win.TopLeftCell = $ws.Range("A22")
$ws.Range("E32").Select() | Out-Null

Write-Output "done"
